# Add a new task-report row (row 11) to the "Report" sheet, matching the
# existing table's layout (Date | Task Name | Status | Person) and reuse
# the date-cell formatting/style used by the other rows in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/number-format (date, centered) from the last existing
# date cell (A10) onto the new date cell (A11) before setting its value,
# so the new row matches the look of the rest of the table.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New row contents:
#   A11 = 4/18/2025 (serial 45765)
#   B11 = "Add exception handling to functions"
#   C11 = "In Progress"
#   D11 = "Adam Rodi"
$ws.Cells.Item(11, 1).Value = 45765
$ws.Cells.Item(11, 2).Value = "Add exception handling to functions"
$ws.Cells.Item(11, 3).Value = "In Progress"
$ws.Cells.Item(11, 4).Value = "Adam Rodi"

# Leave the final selection on D8, matching the saved workbook state.
$null = $ws.Range("D8").Select()
